# Rename 'Codelists' sheet to 'Cells' and update the selected cell on that sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Codelists")
$ws.Name = "Cells"

# Make it the active sheet and move the selection from I17 to G18
$ws.Activate()
$ws.Range("G18").Select()
